$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the default (unstyled) cell style from an untouched data cell
# so we can restore it after forcing number-like strings to stay as text.
$defaultStyle = $ws.Range('D4').Style

$ws.Range('D2').Value = '28.675.58'
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = '1.797.99'
$ws.Range('E3').Value = '  -1.51%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.78'
$ws.Range('D5').Style = $defaultStyle
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5874'
$ws.Range('D6').Style = $defaultStyle
$ws.Range('E6').Value = '  -1.93%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2767'
$ws.Range('D8').Style = $defaultStyle
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('E9').Value = '  -3.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.24'
$ws.Range('D10').Style = $defaultStyle
$ws.Range('E10').Value = '  -0.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07532'
$ws.Range('D11').Style = $defaultStyle
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('D12').Value = '1.789.60'
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.793'
$ws.Range('D13').Style = $defaultStyle
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6201'
$ws.Range('D14').Style = $defaultStyle
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = '2.042.31'
$ws.Range('E15').Value = '  -1.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009115'
$ws.Range('D16').Style = $defaultStyle
$ws.Range('E16').Value = '  -8.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '75.42'
$ws.Range('D17').Style = $defaultStyle
$ws.Range('E17').Value = '  -4.16%  '
$ws.Range('D18').Value = '28.658.50'
$ws.Range('E18').Value = '  -1.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.481'
$ws.Range('D19').Style = $defaultStyle
$ws.Range('E19').Value = '  -5.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.004'
$ws.Range('D20').Style = $defaultStyle
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '211.10'
$ws.Range('D21').Style = $defaultStyle
$ws.Range('E21').Value = '  -5.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.53'
$ws.Range('D22').Style = $defaultStyle
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.834'
$ws.Range('D23').Style = $defaultStyle
$ws.Range('E23').Value = '  -2.13%  '
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.67'
$ws.Range('D25').Style = $defaultStyle
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.965'
$ws.Range('D26').Style = $defaultStyle
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1267'
$ws.Range('D27').Style = $defaultStyle
$ws.Range('E27').Value = '  -1.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.47'
$ws.Range('D28').Style = $defaultStyle
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.423'
$ws.Range('D29').Style = $defaultStyle
$ws.Range('E29').Value = '  -3.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06131'
$ws.Range('D30').Style = $defaultStyle
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.426'
$ws.Range('D31').Style = $defaultStyle
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.818'
$ws.Range('D32').Style = $defaultStyle
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.788'
$ws.Range('D33').Style = $defaultStyle
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.741'
$ws.Range('D34').Style = $defaultStyle
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.051'
$ws.Range('D35').Style = $defaultStyle
$ws.Range('E35').Value = '  -4.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6442'
$ws.Range('D36').Style = $defaultStyle
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.500'
$ws.Range('D37').Style = $defaultStyle
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.719'
$ws.Range('D38').Style = $defaultStyle
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.512'
$ws.Range('D39').Style = $defaultStyle
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01701'
$ws.Range('D40').Style = $defaultStyle
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').Value = '1.144.55'
$ws.Range('E41').Value = '  -6.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8841'
$ws.Range('D42').Style = $defaultStyle
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.007'
$ws.Range('D43').Style = $defaultStyle
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.18'
$ws.Range('D44').Style = $defaultStyle
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('D45').Value = '1.950.09'
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.28'
$ws.Range('D46').Style = $defaultStyle
$ws.Range('E46').Value = '  -3.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000112'
$ws.Range('D47').Style = $defaultStyle
$ws.Range('E47').Value = '  -3.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.601'
$ws.Range('D48').Style = $defaultStyle
$ws.Range('E48').Value = '  +1.98%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05477'
$ws.Range('D49').Style = $defaultStyle
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.336'
$ws.Range('D50').Style = $defaultStyle
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4479'
$ws.Range('D51').Style = $defaultStyle
$ws.Range('E51').Value = '  -1.57%  '
